$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new "2021" data column (I) mirroring the existing 2020 column (H) ---
# Copy H4:H25 formatting onto I4:I25 first so the new column matches the
# existing per-row styling (headers, number formats, borders, etc.)
$ws.Range("H4:H25").Copy()
$ws.Range("I4:I25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header
$ws.Range("I4").Value = 2021

# Data rows (2021 values); rows 6, 9, 13 and 16 are section headers with no
# numeric value in column H either, so they are intentionally left blank.
$ws.Range("I5").Value = 48.5
$ws.Range("I7").Value = 48.8
$ws.Range("I8").Value = 48.2
$ws.Range("I10").Value = 58.2
$ws.Range("I11").Value = 42.4
$ws.Range("I12").Value = 40.7
$ws.Range("I14").Value = 41.5
$ws.Range("I15").Value = 52.6
$ws.Range("I17").Value = 67.1
$ws.Range("I18").Value = 62
$ws.Range("I19").Value = 46.9
$ws.Range("I20").Value = 55.8
$ws.Range("I21").Value = 42.7
$ws.Range("I22").Value = 48.3
$ws.Range("I23").Value = 39.7
$ws.Range("I24").Value = 38.1
$ws.Range("I25").Value = 44.7

# --- Reset the sheet selection back to the top-left cell ---
$ws.Range("A1").Select() | Out-Null
